$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# "Passed" result column (green font) applied to every test-case header row
# ---------------------------------------------------------------------------
$passedRows = @(2,5,9,13,16,19,22,27,31,35,39,43,47,51,56,60,64)
foreach ($r in $passedRows) {
    $cell = $ws.Cells.Item($r, 3)
    $cell.Value = "Passed"
    $cell.Font.Color = 5287936   # RGB(0,176,80) -> FF00B050, BGR-packed long
}

# C10 gets the same style applied but stays empty (style-only touch)
$ws.Cells.Item(10, 3).Font.Color = 5287936

# ---------------------------------------------------------------------------
# Row heights (customHeight) for the existing test-case header rows
# ---------------------------------------------------------------------------
$ws.Rows(2).RowHeight = 29.25
$ws.Rows(5).RowHeight = 39
$ws.Rows(8).RowHeight = 17.25
$ws.Rows(9).RowHeight = 40.5
$ws.Rows(13).RowHeight = 39.75
$ws.Rows(16).RowHeight = 39
$ws.Rows(19).RowHeight = 37.5
$ws.Rows(22).RowHeight = 39.75

# ---------------------------------------------------------------------------
# New Test Case 8 block ("Admin NavBar") appended after the existing cases
# ---------------------------------------------------------------------------
$ws.Range("A51").Value = "Test Case № 8"
$ws.Range("B51").Value = "Admin NavBar"

$ws.Range("A52").Value = "Steps to reproduce :"
$ws.Range("B52").Value = "1. Navigate to the Home Page"
$ws.Range("C52").Value = "Expected Result : "
$ws.Range("D52").Value = "Index Page appears"

$ws.Range("B53").Value = "2. From the NavBar choose ""Registered profiles"" Button"
$ws.Range("C53").Value = "Expected Result : "
$ws.Range("D53").Value = "Registered profiles Page appears"

$ws.Range("A56").Value = "Iteration 2"
$ws.Range("B56").Value = "Check if the NavBar is sending the admin to the chosen page"

$ws.Range("B57").Value = "Do the same but start from Admin Supplements Page"
$ws.Range("C57").Value = "Expected Result : "
$ws.Range("D57").Value = "The admin is redirected to the chosen page"

$ws.Range("A60").Value = "Iteration 3"
$ws.Range("B60").Value = "Check if the NavBar is sending the user to the chosen page"

$ws.Range("B61").Value = "Do the same but start from Admin Fitness Clothing Page"
$ws.Range("C61").Value = "Expected Result : "
$ws.Range("D61").Value = "The admin is redirected to the chosen page"

$ws.Range("A64").Value = "Iteration 4"
$ws.Range("B64").Value = "Check if the NavBar is sending the user to the chosen page"

$ws.Range("B65").Value = "Do the same but start from Admin Messages Page"
$ws.Range("C65").Value = "Expected Result : "
$ws.Range("D65").Value = "The admin is redirected to the chosen page"

# Re-apply "Passed" + green font to the new test case's header/iteration rows
foreach ($r in @(51,56,60,64)) {
    $cell = $ws.Cells.Item($r, 3)
    $cell.Value = "Passed"
    $cell.Font.Color = 5287936
}

# ---------------------------------------------------------------------------
# Sheet view / page setup
# ---------------------------------------------------------------------------
$ws.Range("D9").Select()
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

Write-Output "edit applied"
